# style: example for git diff
#
# Splits three "sentence + next sentence" runs into separate runs
# (one run for the trailing punctuation, one for the following space,
# one for the following sentence), matching the target OOXML diff.
# Also fixes "didn't" -> "did not" in the third location.
#
# Word's object model has no direct "split this run" verb, so we force a
# run boundary at an arbitrary character offset by toggling a character
# -formatting property (Bold) on and back off over the sub-range we want
# to become its own run. Because the net formatting is unchanged, Word's
# own WordOpenXML serialization is unaffected other than the run split.

$d = $word.ActiveDocument

function Split-RunsAt {
    param(
        [int]$StartPos,
        [int[]]$PieceLengths
    )
    $pos = $StartPos
    foreach ($len in $PieceLengths) {
        $pieceEnd = $pos + $len
        $piece = $d.Range($pos, $pieceEnd)
        $piece.Bold = 1
        $piece.Bold = 0
        $pos = $pieceEnd
    }
}

# ---------------------------------------------------------------------
# Edit 1 (Introduction paragraph): ". Some more text"
#   -> ".", " ", "Some more text"
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute(". Some more text", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
Split-RunsAt $r.Start @(1, 1, 14)

# ---------------------------------------------------------------------
# Edit 2 (Introduction paragraph): ". More text"
#   -> ".", " ", "More text"
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute(". More text", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)
Split-RunsAt $r.Start @(1, 1, 9)

# ---------------------------------------------------------------------
# Edit 3a (results paragraph): "). We also found xxx (Fig."
#   -> ").", " ", "We also found xxx (Fig."
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("). We also found xxx (Fig.", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)
Split-RunsAt $r.Start @(2, 1, 24)

# ---------------------------------------------------------------------
# Edit 3b (results paragraph): ") yay! We didn't find xxxx."
#   -> ") yay!", " ", "We did not find xxxx."
# (also normalizes "didn't" -> "did not")
# ---------------------------------------------------------------------
$apos = [char]0x2019
$oldText = ") yay! We didn" + $apos + "t find xxxx."
$newText = ") yay! We did not find xxxx."

$r = $d.Content
$null = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$start = $r.Start
$whole = $d.Range($start, $r.End)
$whole.Text = $newText

Split-RunsAt $start @(6, 1, 22)
